# fixed #99 OneWayRPG-99 スタン状態の実装
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("message")

# Row 33: "HPが足りない" -> "<val1>はスタンした"
$ws.Range("B33").Value = "<val1>はスタンした"

# Row 34: "TPが足りない" -> "<val1>は動けない"
$ws.Range("B34").Value = "<val1>は動けない"

# Row 34 color: yellow -> white
$ws.Range("C34").Value = "white"
